$d = $word.ActiveDocument

# Sanity check: the diary entry we are appending after should be present and
# still end with the "16.9.2022" entry's predecessor text before we touch
# anything.
$fullText = $d.Content.Text
if ($fullText.IndexOf("grid template areas.") -lt 0) {
    throw "Expected anchor text ('grid template areas.') not found in document"
}

# Build the new content: a trailing blank line, an explicit page break, the
# new "16.9.2022" diary entry heading and its paragraph, exactly matching
# the OOXML the author's Word session produced.
$newEntryXml = '<w:p><w:pPr><w:pStyle w:val="BodyText"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>16.9.2022</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>I completed the sixth tutorial</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and make the last two pages left, about my work and how to contact me. I used the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>css</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> grid template columns and flex boxes and learned how to use inheritance in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>css</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and use it to make two kinds of different buttons. </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">I also changed the images used in the tutorial for the </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t>work</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> page by links to some games I made in the past. </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Finally, I completed the last tutorial of the series and deployed my web page in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> pages. I at first had some trouble with some extra folders in my repository and implementing the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> script, but I finally got it all right.</w:t></w:r></w:p>'

# Anchor at the very end of the document's main story and insert the new
# paragraphs there (NB: collapsing any other Range -- e.g. one produced by
# Find.Execute or Paragraphs.Last.Range -- before calling InsertXML on it is
# unreliable in this host and can clobber the preceding paragraph, so we
# deliberately collapse $d.Content itself).
$r = $d.Content
$r.Collapse(0)
$null = $r.InsertXML($newEntryXml)

# Verify the insertion landed correctly.
$finalText = $d.Content.Text
if ($finalText.IndexOf("16.9.2022") -lt 0) {
    throw "Insertion failed: new diary entry not present after InsertXML"
}
if ($finalText.IndexOf("I completed the fourth tutorial") -lt 0) {
    throw "Insertion failed: pre-existing content was lost"
}
